$wb = $excel.ActiveWorkbook

# --- "Relative Samples" sheet: add a "Complex formula" rubric row (A5:D6) ---
$wsRel = $wb.Worksheets.Item("Relative Samples")

$wsRel.Range("A5").Value = "Complex formula"

$wsRel.Range("A6").Formula = "= IF(B6=""ok"",C6,D6)"
$wsRel.Range("B6").Value = "ok"
$wsRel.Range("C6").Value = 100
$wsRel.Range("D6").Value = 200

# New grading rubric comment on A6 (mirrors the A2/A3 rubric comments already on this sheet)
$wsRel.Range("A6").AddComment("rubric:`n score: 1`n type: relative") | Out-Null

# Column A got a bit wider to fit the new content
$wsRel.Range("A1").ColumnWidth = 16.166666666666668

# Update the selection/active cell left behind by editing
$wsRel.Range("C9").Select() | Out-Null

# --- "Relative Samples_CheckOrder" sheet: register the new A6 cell in the grading order ---
$wsRelOrder = $wb.Worksheets.Item("Relative Samples_CheckOrder")
$wsRelOrder.Range("A4").Value = 3
$wsRelOrder.Range("B4").Value = "A6"
$wsRelOrder.Range("G13").Select() | Out-Null

# --- "SheetGradingOrder" sheet: just a leftover selection change ---
$wsOrder = $wb.Worksheets.Item("SheetGradingOrder")
$wsOrder.Range("A5:B5").Select() | Out-Null

# Re-activate "Relative Samples" so it stays the tab shown when the file is opened
$wsRel.Activate() | Out-Null
